$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(5, 1).Value = 42608.901689814818
$ws.Cells.Item(5, 2).Value = -24
$ws.Cells.Item(5, 3).Value = 41
$ws.Cells.Item(5, 4).Value = 55
$ws.Cells.Item(5, 5).Value = 22
$ws.Cells.Item(5, 6).Value = 77
$ws.Cells.Item(5, 7).Value = 24909
$ws.Cells.Item(5, 8).Value = 22752
$ws.Cells.Item(5, 9).Value = 1090
$ws.Cells.Item(5, 10).Value = 216
$ws.Cells.Item(5, 11).Value = 291
$ws.Cells.Item(5, 12).Value = 2
$ws.Cells.Item(5, 13).Value = 7
$ws.Cells.Item(5, 14).Value = "Named"

$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
